# Auto-generated edit script: refresh market-board derived profit columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) for a batch of leves across
# all 8 crafting job sheets, as produced by the scheduled market-data runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 25862.5
$ws.Range("I7").Value = 24900
$ws.Range("J7").Value = 26183.334
$ws.Range("K7").Value = 24900
$ws.Range("L7").Value = 26183.334
$ws.Range("M7").Value = -24788
$ws.Range("N7").Value = -26407.334

$ws.Range("H14").Value = 25862.5
$ws.Range("I14").Value = 24900
$ws.Range("J14").Value = 26183.334
$ws.Range("K14").Value = 24900
$ws.Range("L14").Value = 26183.334
$ws.Range("M14").Value = -24709
$ws.Range("N14").Value = -26565.334

$ws.Range("H80").Value = 9260077
$ws.Range("I80").Value = 20833836
$ws.Range("J80").Value = 1070.6666
$ws.Range("K80").Value = 62501508
$ws.Range("L80").Value = 3211.9998
$ws.Range("M80").Value = -62500510
$ws.Range("N80").Value = -5207.9998

$ws.Range("H83").Value = 9260077
$ws.Range("I83").Value = 20833836
$ws.Range("J83").Value = 1070.6666
$ws.Range("K83").Value = 187504524
$ws.Range("L83").Value = 9635.999400000001
$ws.Range("M83").Value = -187499532
$ws.Range("N83").Value = -19619.9994

$ws.Range("H86").Value = 52634016
$ws.Range("I86").Value = 76925480
$ws.Range("K86").Value = 76925480
$ws.Range("M86").Value = -76924357

$ws.Range("H89").Value = 52634016
$ws.Range("I89").Value = 76925480
$ws.Range("K89").Value = 384627400
$ws.Range("M89").Value = -384621784

$ws.Range("H98").Value = 1205.2
$ws.Range("I98").Value = 1276.6666
$ws.Range("J98").Value = 830
$ws.Range("K98").Value = 1276.6666
$ws.Range("L98").Value = 830
$ws.Range("M98").Value = 221.3334
$ws.Range("N98").Value = -3826

$ws.Range("H111").Value = 4003.4119
$ws.Range("I111").Value = 4117.2
$ws.Range("K111").Value = 12351.6
$ws.Range("M111").Value = -9284.599999999999

$ws.Range("H112").Value = 2618.34
$ws.Range("J112").Value = 2618.34
$ws.Range("L112").Value = 7855.02
$ws.Range("N112").Value = -10071.02

$ws.Range("H113").Value = 141373.2
$ws.Range("I113").Value = 949
$ws.Range("K113").Value = 949
$ws.Range("M113").Value = 2305

$ws.Range("H122").Value = 1205.2
$ws.Range("I122").Value = 1276.6666
$ws.Range("J122").Value = 830
$ws.Range("K122").Value = 3829.9998
$ws.Range("L122").Value = 2490
$ws.Range("M122").Value = -1379.9998
$ws.Range("N122").Value = -7390

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15428.4795
$ws.Range("I32").Value = 14459.985
$ws.Range("K32").Value = 14459.985
$ws.Range("M32").Value = -14172.985

$ws.Range("H61").Value = 5250.081
$ws.Range("I61").Value = 5526.364
$ws.Range("J61").Value = 4844.8667
$ws.Range("K61").Value = 5526.364
$ws.Range("L61").Value = 4844.8667
$ws.Range("M61").Value = -5314.364
$ws.Range("N61").Value = -5268.8667

$ws.Range("H63").Value = 1387.6
$ws.Range("I63").Value = 1609.5
$ws.Range("J63").Value = 500
$ws.Range("K63").Value = 1609.5
$ws.Range("L63").Value = 500
$ws.Range("M63").Value = -923.5
$ws.Range("N63").Value = -1872

$ws.Range("H66").Value = 1387.6
$ws.Range("I66").Value = 1609.5
$ws.Range("J66").Value = 500
$ws.Range("K66").Value = 8047.5
$ws.Range("L66").Value = 2500
$ws.Range("M66").Value = -4615.5
$ws.Range("N66").Value = -9364

$ws.Range("H97").Value = 916.6316
$ws.Range("I97").Value = 617.5
$ws.Range("J97").Value = 1754.2
$ws.Range("K97").Value = 617.5
$ws.Range("L97").Value = 1754.2
$ws.Range("M97").Value = -121.5
$ws.Range("N97").Value = -2746.2

$ws.Range("H132").Value = 2854.3816
$ws.Range("I132").Value = 1170.3684
$ws.Range("K132").Value = 3511.1052
$ws.Range("M132").Value = -981.1052

$ws.Range("H136").Value = 5250.081
$ws.Range("I136").Value = 5526.364
$ws.Range("J136").Value = 4844.8667
$ws.Range("K136").Value = 16579.092
$ws.Range("L136").Value = 14534.6001
$ws.Range("M136").Value = -14029.092
$ws.Range("N136").Value = -19634.6001

$ws.Range("H8").Value = 450
$ws.Range("I8").Value = 450
$ws.Range("K8").Value = 450
$ws.Range("M8").Value = -310

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 320
$ws.Range("I11").Value = 320
$ws.Range("K11").Value = 320
$ws.Range("M11").Value = -180

$ws.Range("H26").Value = 9820
$ws.Range("I26").Value = 9820
$ws.Range("K26").Value = 9820
$ws.Range("M26").Value = -9528

$ws.Range("H105").Value = 2799.0833
$ws.Range("I105").Value = 2850.9524
$ws.Range("J105").Value = 2436
$ws.Range("K105").Value = 2850.9524
$ws.Range("L105").Value = 2436
$ws.Range("M105").Value = -1103.9524
$ws.Range("N105").Value = -5930

$ws.Range("H134").Value = 4714.2036
$ws.Range("I134").Value = 1962.1714
$ws.Range("K134").Value = 5886.5142
$ws.Range("M134").Value = -3351.5142

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 13999.8
$ws.Range("I6").Value = 13999.8
$ws.Range("K6").Value = 13999.8
$ws.Range("M6").Value = -13886.8

$ws.Range("H12").Value = 459.5
$ws.Range("J12").Value = 500
$ws.Range("L12").Value = 500
$ws.Range("N12").Value = -840

$ws.Range("H31").Value = 34485184
$ws.Range("I31").Value = 38463360
$ws.Range("K31").Value = 38463360
$ws.Range("M31").Value = -38463065

$ws.Range("H34").Value = 34485184
$ws.Range("I34").Value = 38463360
$ws.Range("K34").Value = 38463360
$ws.Range("M34").Value = -38463158

$ws.Range("H105").Value = 515
$ws.Range("I105").Value = 515
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 515
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 1232
$ws.Range("N105").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 2396.0588
$ws.Range("I129").Value = 2111.8
$ws.Range("J129").Value = 2514.5
$ws.Range("K129").Value = 6335.400000000001
$ws.Range("L129").Value = 7543.5
$ws.Range("M129").Value = -1335.400000000001
$ws.Range("N129").Value = -17543.5

$ws.Range("H139").Value = 2390.5356
$ws.Range("I139").Value = 1493.45
$ws.Range("K139").Value = 4480.35
$ws.Range("M139").Value = 659.6499999999996

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5482.353
$ws.Range("I70").Value = 4971.2
$ws.Range("K70").Value = 4971.2
$ws.Range("M70").Value = -4701.2

$ws.Range("H73").Value = 5482.353
$ws.Range("I73").Value = 4971.2
$ws.Range("K73").Value = 4971.2
$ws.Range("M73").Value = -4035.2

$ws.Range("H80").Value = 104119.37
$ws.Range("I80").Value = 127886.375
$ws.Range("K80").Value = 127886.375
$ws.Range("M80").Value = -126888.375

$ws.Range("H83").Value = 104119.37
$ws.Range("I83").Value = 127886.375
$ws.Range("K83").Value = 639431.875
$ws.Range("M83").Value = -634439.875

$ws.Range("H123").Value = 40430.9
$ws.Range("J123").Value = 40430.9
$ws.Range("L123").Value = 40430.9
$ws.Range("N123").Value = -45330.9

$ws.Range("H141").Value = 61809.332
$ws.Range("J141").Value = 61809.332
$ws.Range("L141").Value = 61809.332
$ws.Range("N141").Value = -72169.33199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 6066.5713
$ws.Range("J68").Value = 9999.5
$ws.Range("L68").Value = 9999.5
$ws.Range("N68").Value = -11497.5

$ws.Range("H71").Value = 6066.5713
$ws.Range("J71").Value = 9999.5
$ws.Range("L71").Value = 49997.5
$ws.Range("N71").Value = -57485.5

$ws.Range("H82").Value = 3097.1
$ws.Range("I82").Value = 1380.3334
$ws.Range("J82").Value = 5672.25
$ws.Range("K82").Value = 1380.3334
$ws.Range("L82").Value = 5672.25
$ws.Range("M82").Value = -1019.3334
$ws.Range("N82").Value = -6394.25

$ws.Range("H85").Value = 3097.1
$ws.Range("I85").Value = 1380.3334
$ws.Range("J85").Value = 5672.25
$ws.Range("K85").Value = 1380.3334
$ws.Range("L85").Value = 5672.25
$ws.Range("M85").Value = -132.3334
$ws.Range("N85").Value = -8168.25

$ws.Range("H121").Value = 1000
$ws.Range("I121").Value = 1000
$ws.Range("K121").Value = 1000
$ws.Range("M121").Value = 747

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 49999.832
$ws.Range("I2").Value = 49999.832
$ws.Range("K2").Value = 49999.832
$ws.Range("M2").Value = -49887.832

$ws.Range("H4").Value = 575.125
$ws.Range("I4").Value = 514.4286
$ws.Range("K4").Value = 514.4286
$ws.Range("M4").Value = -401.4286

$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

$ws.Range("H107").Value = 906.8461
$ws.Range("I107").Value = 798.7143
$ws.Range("J107").Value = 1033
$ws.Range("K107").Value = 2396.1429
$ws.Range("L107").Value = 3099
$ws.Range("M107").Value = -476.1428999999998
$ws.Range("N107").Value = -6939

$ws.Range("H140").Value = 98990
$ws.Range("J140").Value = 98990
$ws.Range("L140").Value = 98990
$ws.Range("N140").Value = -109350
